$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert four new columns before column E ---
# Old layout: A,B,C,D,E,F,G,H,I
# New layout: A,B,C,D,[E,F,G,H new],I(was E),J(was F),K(was G),L(was H),M(was I)
$ws.Range("E1:H1").EntireColumn.Insert()

# --- New header labels (row 3) for the inserted columns ---
$ws.Range("E3").Value = "客戶編碼"
$ws.Range("F3").Value = "工廠代碼"
$ws.Range("G3").Value = "顏色"
$ws.Range("H3").Value = "尺寸"

# --- Column widths for the 4 new columns ---
# Target raw widths: E,F = 9.875 ; G = 7.875 ; H = 7.25
# ColumnWidth is expressed in "characters" and gets pixel-snapped by the
# engine (MDW=7), so we pick the nearest representable value.
$ws.Columns.Item(5).ColumnWidth = 9.875 - 5/7
$ws.Columns.Item(6).ColumnWidth = 9.875 - 5/7
$ws.Columns.Item(7).ColumnWidth = 7.875 - 5/7
$ws.Columns.Item(8).ColumnWidth = 7.25 - 5/7

# --- Selection moves to E4 ---
$ws.Range("E4").Select()

$wb.Save()
